# Auto-generated edit script applying cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.874.02"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.94%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.806.10"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.90%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "310.09"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.29%  "
$ws.Range("E6").Value = "  +0.05%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4607"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.60%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3733"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.47%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07379"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.59%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8753"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.39"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -3.30%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.351"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "6.525"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "1.724.87"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -5.55%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.07051"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "90.74"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.42%  "
$ws.Range("E17").Value = "  +0.10%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.000008738"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  +0.08%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "14.73"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "26.877.06"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.315"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.51%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.947.99"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -5.21%  "
$ws.Range("E25").Value = "  -2.37%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "151.20"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.27%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "18.39"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.97%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.150"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -9.21%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.290"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.52%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "115.77"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.99%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08888"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.7683"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("E33").Value = "  -3.60%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.474"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  +0.05%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.113"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01956"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05245"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.61%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.406"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.37%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.5341"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "7.234"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.898"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.48%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1660"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.52%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "8.561"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.5055"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.91%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "10.36"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.39%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.02%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "103.81"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.84%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.662"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").Value = "  -1.07%  "
